# Update the "Förändrad" (changed) date in column C for rows 2-146
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C146").Value = 45175
